# Bugfixed evaluation and simulated rt_data for components:
# Column A previously held text labels like "1987Q4" (as shared strings).
# Replace them with real dates (Dec-31 of each year, stored as Excel date
# serial numbers) formatted as "YYYY-MM-DD HH:MM:SS", leaving the header
# row (A1/B1) and column B values untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel date serial numbers (1899-12-30 epoch) for 12/31 of each year
# 1987 .. 2024, corresponding to rows 2 .. 39.
$serials = @(
    32142, 32508, 32873, 33238, 33603, 33969, 34334, 34699, 35064, 35430,
    35795, 36160, 36525, 36891, 37256, 37621, 37986, 38352, 38717, 39082,
    39447, 39813, 40178, 40543, 40908, 41274, 41639, 42004, 42369, 42735,
    43100, 43465, 43830, 44196, 44561, 44926, 45291, 45657
)

for ($i = 0; $i -lt $serials.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $serials[$i]
}

# Apply the date/time number format to the whole updated column range.
$ws.Range("A2:A39").NumberFormat = "YYYY-MM-DD HH:MM:SS"
